$d = $word.ActiveDocument

function Split-ReplaceText($contextText, $narrowText, $replacement) {
    $outer = $d.Content
    $outer.Find.Execute($contextText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $sub = $d.Range($outer.Start, $outer.End)
    $sub.Find.Execute($narrowText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    # Temporarily tweak formatting so the engine keeps this span as an
    # isolated run even after we restore the original formatting value.
    $orig = $sub.Font.Color
    $sub.Font.Color = 255
    $sub.Text = $replacement
    $sub.Font.Color = $orig
}

# 1) "...the parallelisable part takes up 40%..." -> "...the parallelisable section takes up 40%..."
Split-ReplaceText "the parallelisable part takes up 40%" "part " "section "

# 2) ". If the non-parallelisable part of the program takes up 50%..." -> "...non-parallelisable section of the program takes up 50%..."
Split-ReplaceText ". If the non-parallelisable part of the program takes up 50%" "part " "section "

# 3) ". If the parallelisable part of the program takes up 35%..." -> "...parallelisable section of the program takes up 35%..."
Split-ReplaceText ". If the parallelisable part of the program takes up 35%" "part" "section"

# 4a) "If a program takes 0.5 seconds to run" -> "A program takes 0.5 seconds to run"
Split-ReplaceText "If a program takes 0.5 seconds to run" "If a " "A "

# 4b) "...the parallelisable part of the program takes up 50%..." -> "...the parallelisable section of the program takes up 50%..."
Split-ReplaceText "the parallelisable part of the program takes up 50%" "part " "section "
